$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.01632183954996
$ws.Range("C2").Value = 0.4023430362566955
$ws.Range("E2").Value = 0.2104312468631395
$ws.Range("F2").Value = 2.920020281228375
$ws.Range("G2").Value = 0.002481576078569295
$ws.Range("I2").Value = 1.19852442658032
$ws.Range("J2").Value = 0.1024042581127063
$ws.Range("M2").Value = 0.4881751582873264

$ws.Range("B3").Value = 0.9247341823006536
$ws.Range("C3").Value = 0.362348027856342
$ws.Range("E3").Value = 0.2111212489466503
$ws.Range("F3").Value = 2.870403107836196
$ws.Range("G3").Value = 0.002487224107565755
$ws.Range("I3").Value = 1.175547328409138
$ws.Range("J3").Value = 0.1015824560539897
$ws.Range("M3").Value = 0.4681995966316421

$ws.Range("B4").Value = 0.8690486573294152
$ws.Range("C4").Value = 0.3379977774332872
$ws.Range("E4").Value = 0.2115989946675594
$ws.Range("F4").Value = 2.841546261029691
$ws.Range("G4").Value = 0.002490870746426003
$ws.Range("I4").Value = 1.162101214449436
$ws.Range("J4").Value = 0.1011258823149248
$ws.Range("M4").Value = 0.4562506741963333

$ws.Range("B5").Value = 0.8464933194668447
$ws.Range("C5").Value = 0.3281258047595657
$ws.Range("E5").Value = 0.2118073550078741
$ws.Range("F5").Value = 2.830189068588595
$ws.Range("G5").Value = 0.002492401885210676
$ws.Range("I5").Value = 1.156786489278019
$ws.Range("J5").Value = 0.1009517613891404
$ws.Range("M5").Value = 0.4514608657739601

$ws.Range("B6").Value = 0.8427562594529263
$ws.Range("C6").Value = 0.3264896179052812
$ws.Range("E6").Value = 0.211842781269544
$ws.Range("F6").Value = 2.828327443421657
$ws.Range("G6").Value = 0.002492658858425654
$ws.Range("I6").Value = 1.155913877550731
$ws.Range("J6").Value = 0.1009235656682748
$ws.Range("M6").Value = 0.4506703213381016

$ws.Range("B7").Value = 0.8687439152159868
$ws.Range("C7").Value = 0.3378644356385792
$ws.Range("E7").Value = 0.2116017492151556
$ws.Range("F7").Value = 2.841391468386391
$ws.Range("G7").Value = 0.002490891213074196
$ws.Range("I7").Value = 1.162028873797425
$ws.Range("J7").Value = 0.1011234859104384
$ws.Range("M7").Value = 0.4561857554015774

$ws.Range("B8").Value = 0.9846275367516455
$ws.Range("C8").Value = 0.3885091711360928
$ws.Range("E8").Value = 0.2106579800146573
$ws.Range("F8").Value = 2.902577119732896
$ws.Range("G8").Value = 0.002483486521931261
$ws.Range("I8").Value = 1.190463467965287
$ws.Range("J8").Value = 0.102110851872105
$ws.Range("M8").Value = 0.4812219408200846

$ws.Range("B9").Value = 1.216306982365381
$ws.Range("C9").Value = 0.4895180781534236
$ws.Range("E9").Value = 0.2092331113738339
$ws.Range("F9").Value = 3.035432585774146
$ws.Range("G9").Value = 0.002470376652586503
$ws.Range("I9").Value = 1.251558989414718
$ws.Range("J9").Value = 0.1044344072454564
$ws.Range("M9").Value = 0.5328313831683573

$ws.Range("B10").Value = 1.38934205018495
$ws.Range("C10").Value = 0.5648447060651733
$ws.Range("E10").Value = 0.2084416135365466
$ws.Range("F10").Value = 3.141053301013528
$ws.Range("G10").Value = 0.002461594440167057
$ws.Range("I10").Value = 1.299819012352884
$ws.Range("J10").Value = 0.1063866571600443
$ws.Range("M10").Value = 0.5722937856243107

$ws.Range("B11").Value = 1.468700008489236
$ws.Range("C11").Value = 0.59937343004799
$ws.Range("E11").Value = 0.2081361081262649
$ws.Range("F11").Value = 3.190879490617448
$ws.Range("G11").Value = 0.002457781442831554
$ws.Range("I11").Value = 1.322532601180939
$ws.Range("J11").Value = 0.1073299780978374
$ws.Range("M11").Value = 0.5905849989595708

$ws.Range("B12").Value = 1.498845231437997
$ws.Range("C12").Value = 0.6124876423402839
$ws.Range("E12").Value = 0.2080281938499606
$ws.Range("F12").Value = 3.210005924578013
$ws.Range("G12").Value = 0.002456363571631375
$ws.Range("I12").Value = 1.331244965673932
$ws.Range("J12").Value = 0.1076952881290723
$ws.Range("M12").Value = 0.5975604167660009

$ws.Range("B13").Value = 1.492348708374095
$ws.Range("C13").Value = 0.6096615108016863
$ws.Range("E13").Value = 0.2080510904788504
$ws.Range("F13").Value = 3.205875174251446
$ws.Range("G13").Value = 0.00245666778023882
$ws.Range("I13").Value = 1.329363624224627
$ws.Range("J13").Value = 0.1076162496165622
$ws.Range("M13").Value = 0.5960559576742241

$ws.Range("B14").Value = 1.47117817936379
$ws.Range("C14").Value = 0.6004515568168358
$ws.Range("E14").Value = 0.2081270745504256
$ws.Range("F14").Value = 3.192447841674039
$ws.Range("G14").Value = 0.002457664272934799
$ws.Range("I14").Value = 1.323247132367101
$ws.Range("J14").Value = 0.1073598692632558
$ws.Range("M14").Value = 0.5911578894667571

$ws.Range("B15").Value = 1.458222919112359
$ws.Range("C15").Value = 0.594815301039489
$ws.Range("E15").Value = 0.2081746273225846
$ws.Range("F15").Value = 3.184256932095309
$ws.Range("G15").Value = 0.002458278038614609
$ws.Range("I15").Value = 1.31951514703178
$ws.Range("J15").Value = 0.1072038876293462
$ws.Range("M15").Value = 0.5881640569754438

$ws.Range("B16").Value = 1.384168900824761
$ws.Range("C16").Value = 0.5625935539852662
$ws.Range("E16").Value = 0.2084626706972887
$ws.Range("F16").Value = 3.137833053402119
$ws.Range("G16").Value = 0.002461847279052217
$ws.Range("I16").Value = 1.298350082998624
$ws.Range("J16").Value = 0.1063261333490573
$ws.Range("M16").Value = 0.5711052532703889

$ws.Range("B17").Value = 1.338905053154519
$ws.Range("C17").Value = 0.5428945890622572
$ws.Range("E17").Value = 0.208653297356804
$ws.Range("F17").Value = 3.109810903187594
$ws.Range("G17").Value = 0.002464083416566884
$ws.Range("I17").Value = 1.285562006683051
$ws.Range("J17").Value = 0.105801914076288
$ws.Range("M17").Value = 0.560727270109652

$ws.Range("B18").Value = 1.312930952038528
$ws.Range("C18").Value = 0.5315888765042587
$ws.Range("E18").Value = 0.2087680794972009
$ws.Range("F18").Value = 3.093860578025669
$ws.Range("G18").Value = 0.002465386730559804
$ws.Range("I18").Value = 1.278278047622891
$ws.Range("J18").Value = 0.1055055865945462
$ws.Range("M18").Value = 0.5547900834806612

$ws.Range("B19").Value = 1.304146906239907
$ws.Range("C19").Value = 0.5277651466451516
$ws.Range("E19").Value = 0.2088078275018379
$ws.Range("F19").Value = 3.088488724388895
$ws.Range("G19").Value = 0.002465830960189562
$ws.Range("I19").Value = 1.275824028658107
$ws.Range("J19").Value = 0.1054061418876628
$ws.Range("M19").Value = 0.5527853394636395

$ws.Range("B20").Value = 1.343717197778346
$ws.Range("C20").Value = 0.544989021587071
$ws.Range("E20").Value = 0.2086324734979925
$ws.Range("F20").Value = 3.112776578287253
$ws.Range("G20").Value = 0.002463843602407595
$ws.Range("I20").Value = 1.286915915065649
$ws.Range("J20").Value = 0.1058571799105152
$ws.Range("M20").Value = 0.5618287162337481

$ws.Range("B21").Value = 1.477393910736225
$ws.Range("C21").Value = 0.6031556770559519
$ws.Range("E21").Value = 0.2081045457580348
$ws.Range("F21").Value = 3.1963847430616
$ws.Range("G21").Value = 0.002457370873785009
$ws.Range("I21").Value = 1.325040660131521
$ws.Range("J21").Value = 0.1074349534747441
$ws.Range("M21").Value = 0.5925952419342337

$ws.Range("B22").Value = 1.565308539951843
$ws.Range("C22").Value = 0.6413984685966057
$ws.Range("E22").Value = 0.2078047953809055
$ws.Range("F22").Value = 3.252534505616069
$ws.Range("G22").Value = 0.002453292211494708
$ws.Range("I22").Value = 1.350606585232327
$ws.Range("J22").Value = 0.1085133743586937
$ws.Range("M22").Value = 0.6129882174343066

$ws.Range("B23").Value = 1.518336069290797
$ws.Range("C23").Value = 0.6209663435492985
$ws.Range("E23").Value = 0.207960657549302
$ws.Range("F23").Value = 3.222427580011868
$ws.Range("G23").Value = 0.002455455245983697
$ws.Range("I23").Value = 1.336901518592256
$ws.Range("J23").Value = 0.1079334263545348
$ws.Range("M23").Value = 0.6020779626071402

$ws.Range("B24").Value = 1.341541475658914
$ws.Range("C24").Value = 0.5440420680433249
$ws.Range("E24").Value = 0.2086418717987719
$ws.Range("F24").Value = 3.111435298457764
$ws.Range("G24").Value = 0.00246395196714593
$ws.Range("I24").Value = 1.286303601190852
$ws.Range("J24").Value = 0.1058321784987299
$ws.Range("M24").Value = 0.5613306613284337

$ws.Range("B25").Value = 1.153145029796178
$ws.Range("C25").Value = 0.4620026376072701
$ws.Range("E25").Value = 0.2095734296736325
$ws.Range("F25").Value = 2.998096768494833
$ws.Range("G25").Value = 0.002473773270050779
$ws.Range("I25").Value = 1.234447424482951
$ws.Range("J25").Value = 0.1037634297890762
$ws.Range("M25").Value = 0.5185993133476501

